$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R2").Value = "'712441"
$ws.Range("R2").Style = "Normal"
$ws.Range("S2").Value = "'12874"
$ws.Range("S2").Style = "Normal"
$ws.Range("R3").Value = "'19485"
$ws.Range("R3").Style = "Normal"
$ws.Range("R4").Value = "'1203674"
$ws.Range("R4").Style = "Normal"
$ws.Range("S4").Value = "'1"
$ws.Range("S4").Style = "Normal"
$ws.Range("R5").Value = "'1406805"
$ws.Range("R5").Style = "Normal"
$ws.Range("S5").Value = "'182347"
$ws.Range("S5").Style = "Normal"
$ws.Range("R6").Value = "'81439"
$ws.Range("R6").Style = "Normal"
$ws.Range("S6").Value = "'12215"
$ws.Range("S6").Style = "Normal"
$ws.Range("R7").Value = "'102147"
$ws.Range("R7").Style = "Normal"
$ws.Range("S7").Value = "'21360"
$ws.Range("S7").Style = "Normal"
$ws.Range("R8").Value = "'37735516"
$ws.Range("R8").Style = "Normal"
$ws.Range("S8").Value = "'579076"
$ws.Range("S8").Style = "Normal"
$ws.Range("R9").Value = "'25872729"
$ws.Range("R9").Style = "Normal"
$ws.Range("S9").Value = "'23749"
$ws.Range("S9").Style = "Normal"
$ws.Range("R10").Value = "'10176957"
$ws.Range("R10").Style = "Normal"
$ws.Range("S10").Value = "'13413163"
$ws.Range("S10").Style = "Normal"
$ws.Range("R11").Value = "'60781125"
$ws.Range("R11").Style = "Normal"
$ws.Range("S11").Value = "'97390"
$ws.Range("S11").Style = "Normal"
$ws.Range("R12").Value = "'2963807"
$ws.Range("R12").Style = "Normal"
$ws.Range("S12").Value = "'360191"
$ws.Range("S12").Style = "Normal"
$ws.Range("R13").Value = "'29062706"
$ws.Range("R13").Style = "Normal"
$ws.Range("S13").Value = "'2069694"
$ws.Range("S13").Style = "Normal"
$ws.Range("R14").Value = "'658649"
$ws.Range("R14").Style = "Normal"
$ws.Range("S14").Value = "'197056"
$ws.Range("S14").Style = "Normal"
$ws.Range("R15").Value = "'102180"
$ws.Range("R15").Style = "Normal"
$ws.Range("S15").Value = "'4044660"
$ws.Range("S15").Style = "Normal"
$ws.Range("R16").Value = "'156149810"
$ws.Range("R16").Style = "Normal"
$ws.Range("S16").Value = "'122631498"
$ws.Range("S16").Style = "Normal"
$ws.Range("R17").Value = "'158472"
$ws.Range("R17").Style = "Normal"
$ws.Range("R18").Value = "'4879642"
$ws.Range("R18").Style = "Normal"
$ws.Range("S18").Value = "'632628"
$ws.Range("S18").Style = "Normal"
$ws.Range("R19").Value = "'2745065"
$ws.Range("R19").Style = "Normal"
$ws.Range("S19").Value = "'885"
$ws.Range("S19").Style = "Normal"
$ws.Range("R20").Value = "'7684727"
$ws.Range("R20").Style = "Normal"
$ws.Range("S20").Value = "'2907579"
$ws.Range("S20").Style = "Normal"
$ws.Range("R21").Value = "'2384371"
$ws.Range("R21").Style = "Normal"
$ws.Range("S21").Value = "'7742765"
$ws.Range("S21").Style = "Normal"
$ws.Range("R22").Value = "'15801866"
$ws.Range("R22").Style = "Normal"
$ws.Range("S22").Value = "'1500146"
$ws.Range("S22").Style = "Normal"
$ws.Range("R23").Value = "'4217016"
$ws.Range("R23").Style = "Normal"
$ws.Range("S23").Value = "'396421"
$ws.Range("S23").Style = "Normal"
$ws.Range("R24").Value = "'25547563"
$ws.Range("R24").Style = "Normal"
$ws.Range("S24").Value = "'8431030"
$ws.Range("S24").Style = "Normal"
$ws.Range("R25").Value = "'4417026"
$ws.Range("R25").Style = "Normal"
$ws.Range("S25").Value = "'234570"
$ws.Range("S25").Style = "Normal"
$ws.Range("R26").Value = "'11500429"
$ws.Range("R26").Style = "Normal"
$ws.Range("S26").Value = "'3005434"
$ws.Range("S26").Style = "Normal"
$ws.Range("R27").Value = "'2796217"
$ws.Range("R27").Style = "Normal"
$ws.Range("S27").Value = "'164984"
$ws.Range("S27").Style = "Normal"
$ws.Range("R28").Value = "'326144631"
$ws.Range("R28").Style = "Normal"
$ws.Range("S28").Value = "'24807"
$ws.Range("S28").Style = "Normal"
$ws.Range("R29").Value = "'5970148"
$ws.Range("R29").Style = "Normal"
$ws.Range("R30").Value = "'16591276"
$ws.Range("R30").Style = "Normal"
$ws.Range("S30").Value = "'22164"
$ws.Range("S30").Style = "Normal"
$ws.Range("R31").Value = "'42469335"
$ws.Range("R31").Style = "Normal"
$ws.Range("S31").Value = "'2109946"
$ws.Range("S31").Style = "Normal"
$ws.Range("R32").Value = "'31204503"
$ws.Range("R32").Style = "Normal"
$ws.Range("R33").Value = "'10526266"
$ws.Range("R33").Style = "Normal"
$ws.Range("S33").Value = "'422001"
$ws.Range("S33").Style = "Normal"
$ws.Range("R34").Value = "'18660423"
$ws.Range("R34").Style = "Normal"
$ws.Range("S34").Value = "'1600041"
$ws.Range("S34").Style = "Normal"
$ws.Range("R35").Value = "'9365056"
$ws.Range("R35").Style = "Normal"
$ws.Range("S35").Value = "'11578"
$ws.Range("S35").Style = "Normal"
$ws.Range("R36").Value = "'3496904"
$ws.Range("R36").Style = "Normal"
$ws.Range("S36").Value = "'3274"
$ws.Range("S36").Style = "Normal"
$ws.Range("R37").Value = "'1580848"
$ws.Range("R37").Style = "Normal"
$ws.Range("R38").Value = "'1397953"
$ws.Range("R38").Style = "Normal"
$ws.Range("S38").Value = "'14"
$ws.Range("S38").Style = "Normal"
$ws.Range("R39").Value = "'16985407"
$ws.Range("R39").Style = "Normal"
$ws.Range("S39").Value = "'3870645"
$ws.Range("S39").Style = "Normal"
$ws.Range("R40").Value = "'67314905"
$ws.Range("R40").Style = "Normal"
$ws.Range("S40").Value = "'1560201"
$ws.Range("S40").Style = "Normal"
$ws.Range("R41").Value = "'18491717"
$ws.Range("R41").Style = "Normal"
$ws.Range("R42").Value = "'81618"
$ws.Range("R42").Style = "Normal"
$ws.Range("S42").Value = "'479632"
$ws.Range("S42").Style = "Normal"
$ws.Range("R43").Value = "'2853773"
$ws.Range("R43").Style = "Normal"
$ws.Range("S43").Value = "'264914"
$ws.Range("S43").Style = "Normal"
$ws.Range("R45").Value = "'4817370"
$ws.Range("R45").Style = "Normal"
$ws.Range("S45").Value = "'9428698"
$ws.Range("S45").Style = "Normal"
$ws.Range("R46").Value = "'3625"
$ws.Range("R46").Style = "Normal"
$ws.Range("R47").Value = "'10485"
$ws.Range("R47").Style = "Normal"
$ws.Range("S47").Value = "'8428"
$ws.Range("S47").Style = "Normal"
$ws.Range("R48").Value = "'458607"
$ws.Range("R48").Style = "Normal"
$ws.Range("R49").Value = "'17442459"
$ws.Range("R49").Style = "Normal"
$ws.Range("S49").Value = "'1222664"
$ws.Range("S49").Style = "Normal"
$ws.Range("R50").Value = "'5526568"
$ws.Range("R50").Style = "Normal"
$ws.Range("S50").Value = "'19815"
$ws.Range("S50").Style = "Normal"
$ws.Range("R51").Value = "'757102"
$ws.Range("R51").Style = "Normal"
$ws.Range("R52").Value = "'5144822"
$ws.Range("R52").Style = "Normal"
$ws.Range("R53").Value = "'11840829"
$ws.Range("R53").Style = "Normal"
$ws.Range("S53").Value = "'29754"
$ws.Range("S53").Style = "Normal"
$ws.Range("R54").Value = "'5417336"
$ws.Range("R54").Style = "Normal"
$ws.Range("S54").Value = "'6661427"
$ws.Range("S54").Style = "Normal"
$ws.Range("R55").Value = "'3860627"
$ws.Range("R55").Style = "Normal"
$ws.Range("S55").Value = "'3133411"
$ws.Range("S55").Style = "Normal"
$ws.Range("R56").Value = "'20741854"
$ws.Range("R56").Style = "Normal"
$ws.Range("S56").Value = "'14321208"
$ws.Range("S56").Style = "Normal"
$ws.Range("R57").Value = "'1585724"
$ws.Range("R57").Style = "Normal"
$ws.Range("S57").Value = "'5223204"
$ws.Range("S57").Style = "Normal"
$ws.Range("R58").Value = "'577804"
$ws.Range("R58").Style = "Normal"
$ws.Range("S58").Value = "'12253236"
$ws.Range("S58").Style = "Normal"
$ws.Range("R59").Value = "'1740805"
$ws.Range("R59").Style = "Normal"
$ws.Range("S59").Value = "'6227"
$ws.Range("S59").Style = "Normal"
$ws.Range("R60").Value = "'3169805"
$ws.Range("R60").Style = "Normal"
$ws.Range("S60").Value = "'139654"
$ws.Range("S60").Style = "Normal"
$ws.Range("R61").Value = "'10405980"
$ws.Range("R61").Style = "Normal"
$ws.Range("S61").Value = "'1370"
$ws.Range("S61").Style = "Normal"
$ws.Range("R62").Value = "'12111750"
$ws.Range("R62").Style = "Normal"
$ws.Range("S62").Value = "'4177835"
$ws.Range("S62").Style = "Normal"
$ws.Range("R63").Value = "'12111750"
$ws.Range("R63").Style = "Normal"
$ws.Range("S63").Value = "'4177835"
$ws.Range("S63").Style = "Normal"
$ws.Range("R64").Value = "'2708165"
$ws.Range("R64").Style = "Normal"
$ws.Range("S64").Value = "'2879529"
$ws.Range("S64").Style = "Normal"
$ws.Range("R65").Value = "'9785731"
$ws.Range("R65").Style = "Normal"
$ws.Range("S65").Value = "'2043732"
$ws.Range("S65").Style = "Normal"
$ws.Range("R66").Value = "'1444216"
$ws.Range("R66").Style = "Normal"
$ws.Range("S66").Value = "'72112"
$ws.Range("S66").Style = "Normal"
$ws.Range("R67").Value = "'790075"
$ws.Range("R67").Style = "Normal"
$ws.Range("S67").Value = "'16"
$ws.Range("S67").Style = "Normal"
$ws.Range("R68").Value = "'164650"
$ws.Range("R68").Style = "Normal"
$ws.Range("S68").Value = "'95814"
$ws.Range("S68").Style = "Normal"
$ws.Range("R69").Value = "'5423689"
$ws.Range("R69").Style = "Normal"
$ws.Range("S69").Value = "'101102"
$ws.Range("S69").Style = "Normal"
$ws.Range("R70").Value = "'7754019"
$ws.Range("R70").Style = "Normal"
$ws.Range("S70").Value = "'123100"
$ws.Range("S70").Style = "Normal"
$ws.Range("R71").Value = "'10559952"
$ws.Range("R71").Style = "Normal"
$ws.Range("S71").Value = "'990118"
$ws.Range("S71").Style = "Normal"
$ws.Range("R72").Value = "'32572850"
$ws.Range("R72").Style = "Normal"
$ws.Range("S72").Value = "'449586"
$ws.Range("S72").Style = "Normal"
$ws.Range("R73").Value = "'136890453"
$ws.Range("R73").Style = "Normal"
$ws.Range("S73").Value = "'14304275"
$ws.Range("S73").Style = "Normal"
$ws.Range("R74").Value = "'25608357"
$ws.Range("R74").Style = "Normal"
$ws.Range("S74").Value = "'2052917"
$ws.Range("S74").Style = "Normal"
$ws.Range("R75").Value = "'9949912"
$ws.Range("R75").Style = "Normal"
$ws.Range("S75").Value = "'1329969"
$ws.Range("S75").Style = "Normal"
$ws.Range("R76").Value = "'23590"
$ws.Range("R76").Style = "Normal"
$ws.Range("R77").Value = "'15064968"
$ws.Range("R77").Style = "Normal"
$ws.Range("S77").Value = "'271364"
$ws.Range("S77").Style = "Normal"
$ws.Range("R78").Value = "'84221"
$ws.Range("R78").Style = "Normal"
$ws.Range("S78").Value = "'937631"
$ws.Range("S78").Style = "Normal"
$ws.Range("R79").Value = "'4978820"
$ws.Range("R79").Style = "Normal"
$ws.Range("S79").Value = "'37209"
$ws.Range("S79").Style = "Normal"
$ws.Range("R80").Value = "'386227"
$ws.Range("R80").Style = "Normal"
$ws.Range("R81").Value = "'50106"
$ws.Range("R81").Style = "Normal"
$ws.Range("R82").Value = "'5106154"
$ws.Range("R82").Style = "Normal"
$ws.Range("S82").Value = "'2284419"
$ws.Range("S82").Style = "Normal"
$ws.Range("R83").Value = "'6544932"
$ws.Range("R83").Style = "Normal"
$ws.Range("S83").Value = "'417227"
$ws.Range("S83").Style = "Normal"
$ws.Range("R84").Value = "'124154272"
$ws.Range("R84").Style = "Normal"
$ws.Range("S84").Value = "'1531198"
$ws.Range("S84").Style = "Normal"
$ws.Range("R85").Value = "'116957970"
$ws.Range("R85").Style = "Normal"
$ws.Range("S85").Value = "'127834"
$ws.Range("S85").Style = "Normal"
$ws.Range("R86").Value = "'30475"
$ws.Range("R86").Style = "Normal"
$ws.Range("R87").Value = "'109121119"
$ws.Range("R87").Style = "Normal"
$ws.Range("S87").Value = "'111157"
$ws.Range("S87").Style = "Normal"
$ws.Range("R88").Value = "'7557708"
$ws.Range("R88").Style = "Normal"
$ws.Range("R89").Value = "'23920"
$ws.Range("R89").Style = "Normal"
$ws.Range("R90").Value = "'27573807"
$ws.Range("R90").Style = "Normal"
$ws.Range("S90").Value = "'158453"
$ws.Range("S90").Style = "Normal"
$ws.Range("R91").Value = "'1552361"
$ws.Range("R91").Style = "Normal"
$ws.Range("S91").Value = "'564"
$ws.Range("S91").Style = "Normal"
$ws.Range("R92").Value = "'416460"
$ws.Range("R92").Style = "Normal"
$ws.Range("S92").Value = "'1103747"
$ws.Range("S92").Style = "Normal"
$ws.Range("S93").Value = "'2243"
$ws.Range("S93").Style = "Normal"
$ws.Range("R94").Value = "'8537486"
$ws.Range("R94").Style = "Normal"
$ws.Range("S94").Value = "'168203"
$ws.Range("S94").Style = "Normal"
$ws.Range("R95").Value = "'4035127"
$ws.Range("R95").Style = "Normal"
$ws.Range("S95").Value = "'391062"
$ws.Range("S95").Style = "Normal"
$ws.Range("R96").Value = "'10272844"
$ws.Range("R96").Style = "Normal"
$ws.Range("S96").Value = "'406746"
$ws.Range("S96").Style = "Normal"
$ws.Range("R97").Value = "'4909"
$ws.Range("R97").Style = "Normal"
$ws.Range("S97").Value = "'269283"
$ws.Range("S97").Style = "Normal"
